# Apply the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Wed Jun 28 16:46:58 UTC 2023 with GitHub Actions".
# Rows 2-51 on Sheet1 hold one coin each (B=Coin, C=Link, D=Price, E=Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.338.66"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").Value = "1.858.61"
$ws.Range("E3").Value = "  -1.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").Value = "'234.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

# Row 6
$ws.Range("E6").Value = "  +0.32%  "

# Row 7
$ws.Range("D7").Value = "'0.4732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.14%  "

# Row 8
$ws.Range("D8").Value = "'0.2739"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.24%  "

# Row 9
$ws.Range("D9").Value = "'0.06430"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.58%  "

# Row 10
$ws.Range("D10").Value = "1.853.93"
$ws.Range("E10").Value = "  -6.87%  "

# Row 11
$ws.Range("D11").Value = "'0.07460"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "

# Row 12
$ws.Range("D12").Value = "'16.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.90%  "

# Row 13
$ws.Range("D13").Value = "'4.988"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "

# Row 14
$ws.Range("D14").Value = "'85.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.82%  "

# Row 15
$ws.Range("D15").Value = "'0.6328"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.97%  "

# Row 16
$ws.Range("D16").Value = "30.305.32"
$ws.Range("E16").Value = "  -0.98%  "

# Row 17
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "

# Row 18
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'230.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.10%  "

# Row 20
$ws.Range("D20").Value = "'0.000007414"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "

# Row 21
$ws.Range("D21").Value = "2.099.54"
$ws.Range("E21").Value = "  -5.96%  "

# Row 22
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("D23").Value = "'4.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.59%  "

# Row 24
$ws.Range("D24").Value = "'5.995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.76%  "

# Row 25
$ws.Range("D25").Value = "'9.261"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "

# Row 26
$ws.Range("D26").Value = "'166.38"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'17.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.90%  "

# Row 28
$ws.Range("D28").Value = "'1.890"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$ws.Range("D29").Value = "'0.1047"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.37%  "

# Row 30
$ws.Range("D30").Value = "'1.402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("D31").Value = "'4.148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.61%  "

# Row 32
$ws.Range("D32").Value = "'3.930"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.03%  "

# Row 33
$ws.Range("D33").Value = "'0.04938"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.29%  "

# Row 34
$ws.Range("D34").Value = "'1.164"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.23%  "

# Row 35
$ws.Range("D35").Value = "'0.7246"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.91%  "

# Row 36
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").Value = "'2.702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "

# Row 38
$ws.Range("D38").Value = "'0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39
$ws.Range("D39").Value = "'2.651"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("D40").Value = "'0.9168"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.50%  "

# Row 41
$ws.Range("D41").Value = "'1.970"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.52%  "

# Row 42
$ws.Range("D42").Value = "'106.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

# Row 43
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").Value = "'0.4112"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.02%  "

# Row 45
$ws.Range("D45").Value = "'5.581"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.42%  "

# Row 46
$ws.Range("D46").Value = "'7.113"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "

# Row 47
$ws.Range("D47").Value = "'60.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.86%  "

# Row 48
$ws.Range("D48").Value = "'0.1198"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.75%  "

# Row 49
$ws.Range("D49").Value = "'8.699"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.91%  "

# Row 50
$ws.Range("D50").Value = "'33.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "

# Row 51
$ws.Range("D51").Value = "'1.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.70%  "
